# Apply the "updated cost model with better dbe cost calculations" edit.
#
# The old "SiteDevelopmentValues" sheet had a row for "timing_dbe_cost"
# (labelled "Maser") that bundled the maser and DBE/BDC costs together.
# This change splits that single row into two separate, clearly-labelled
# cost rows:
#   - row 14: maser_cost       = "Cost for maser & components" ($500,000, from Jono)
#   - row 15: dbe_cost (new)   = "Cost for BDC & DBE"           ($250,000, from Ranjani - assumes 2SB, 2pol)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SiteDevelopmentValues")

# Insert a new row below row 14 so the old row 14 becomes "maser_cost" and
# the newly inserted row 15 becomes "dbe_cost"; everything below shifts down.
$ws.Rows.Item(15).Insert()

# Row 14: rename/refresh the maser cost entry and its value/comment.
$ws.Range("A14").Value = "maser_cost"
$ws.Range("B14").Value = "Cost for maser & components"
$ws.Range("C14").Value = 500000
$ws.Range("D14").Value = "from Jono"

# Row 15 (new): the split-out DBE/BDC cost entry.
$ws.Range("A15").Value = "dbe_cost"
$ws.Range("B15").Value = "Cost for BDC & DBE"
$ws.Range("C15").Value = 250000
$ws.Range("D15").Value = "from Ranjani - assumes 2SB, 2pol"

# Reflect the author's final UI state: SiteDevelopmentValues tab active,
# with cell C15 (the new dbe_cost value) selected.
$ws.Activate()
$ws.Range("C15").Select() | Out-Null
